$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44355
$ws.Range("M2").Value = 150
$ws.Range("N2").Value = 16800
$ws.Range("O2").Value = 18200
$ws.Range("P2").Value = 17500
$ws.Range("S2").Value = 1250
$ws.Range("T2").Value = 14
$ws.Range("K2").Value = 'Mankaki'
$ws.Range("L2").Value = 'Primera'
$ws.Range("Q2").Value = '$/caja 14 kilos granel'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("D3").Value = 44355
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 15400
$ws.Range("O3").Value = 15400
$ws.Range("P3").Value = 15400
$ws.Range("S3").Value = 1100
$ws.Range("T3").Value = 14
$ws.Range("K3").Value = 'Mankaki'
$ws.Range("L3").Value = 'Segunda'
$ws.Range("Q3").Value = '$/caja 14 kilos granel'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("D4").Value = 44320
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 10000
$ws.Range("S4").Value = 625
$ws.Range("T4").Value = 16
$ws.Range("K4").Value = 'Mankaki'
$ws.Range("L4").Value = 'Primera'
$ws.Range("Q4").Value = '$/caja 16 kilos granel'
$ws.Range("R4").Value = 'Provincia de Quillota'
$ws.Range("D5").Value = 45062
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 14000
$ws.Range("S5").Value = 875
$ws.Range("T5").Value = 16
$ws.Range("K5").Value = 'Fuyu'
$ws.Range("L5").Value = 'Primera'
$ws.Range("Q5").Value = '$/caja 16 kilos granel'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("D6").Value = 45062
$ws.Range("M6").Value = 470
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 13500
$ws.Range("P6").Value = 13266
$ws.Range("S6").Value = 829
$ws.Range("T6").Value = 16
$ws.Range("K6").Value = 'Mankaki'
$ws.Range("L6").Value = 'Primera'
$ws.Range("Q6").Value = '$/caja 16 kilos granel'
$ws.Range("R6").Value = 'Provincia de Curicó'
$ws.Range("D7").Value = 45072
$ws.Range("M7").Value = 470
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 14000
$ws.Range("P7").Value = 13468
$ws.Range("S7").Value = 842
$ws.Range("T7").Value = 16
$ws.Range("K7").Value = 'Fuyu'
$ws.Range("L7").Value = 'Primera'
$ws.Range("Q7").Value = '$/caja 16 kilos granel'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("D8").Value = 45084
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 11000
$ws.Range("O8").Value = 11000
$ws.Range("P8").Value = 11000
$ws.Range("S8").Value = 688
$ws.Range("T8").Value = 16
$ws.Range("K8").Value = 'Hachiya'
$ws.Range("L8").Value = 'Primera'
$ws.Range("Q8").Value = '$/caja 16 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("D9").Value = 45084
$ws.Range("M9").Value = 220
$ws.Range("N9").Value = 8000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 8000
$ws.Range("S9").Value = 500
$ws.Range("T9").Value = 16
$ws.Range("K9").Value = 'Hachiya'
$ws.Range("L9").Value = 'Segunda'
$ws.Range("Q9").Value = '$/caja 16 kilos granel'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("D10").Value = 44319
$ws.Range("M10").Value = 65
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 15000
$ws.Range("S10").Value = 1000
$ws.Range("T10").Value = 15
$ws.Range("K10").Value = 'Mankaki'
$ws.Range("L10").Value = 'Primera'
$ws.Range("Q10").Value = '$/caja 15 kilos granel'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("D11").Value = 44321
$ws.Range("M11").Value = 95
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 15000
$ws.Range("S11").Value = 1250
$ws.Range("T11").Value = 12
$ws.Range("K11").Value = 'Fuyu'
$ws.Range("L11").Value = 'Especial'
$ws.Range("Q11").Value = '$/caja 12 kilos empedrada'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("D12").Value = 44321
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 13000
$ws.Range("O12").Value = 13000
$ws.Range("P12").Value = 13000
$ws.Range("S12").Value = 1083
$ws.Range("T12").Value = 12
$ws.Range("K12").Value = 'Fuyu'
$ws.Range("L12").Value = 'Primera'
$ws.Range("Q12").Value = '$/caja 12 kilos empedrada'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("D13").Value = 44321
$ws.Range("M13").Value = 4
$ws.Range("N13").Value = 270000
$ws.Range("O13").Value = 270000
$ws.Range("P13").Value = 270000
$ws.Range("S13").Value = 600
$ws.Range("T13").Value = 450
$ws.Range("K13").Value = 'Mankaki'
$ws.Range("L13").Value = 'Primera'
$ws.Range("Q13").Value = '$/bins (450 kilos)'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("D14").Value = 44321
$ws.Range("M14").Value = 5
$ws.Range("N14").Value = 240000
$ws.Range("O14").Value = 240000
$ws.Range("P14").Value = 240000
$ws.Range("S14").Value = 533
$ws.Range("T14").Value = 450
$ws.Range("K14").Value = 'Mankaki'
$ws.Range("L14").Value = 'Segunda'
$ws.Range("Q14").Value = '$/bins (450 kilos)'
$ws.Range("R14").Value = 'Región de O''Higgins'
$ws.Range("D15").Value = 45083
$ws.Range("M15").Value = 150
$ws.Range("N15").Value = 14400
$ws.Range("O15").Value = 14400
$ws.Range("P15").Value = 14400
$ws.Range("S15").Value = 900
$ws.Range("T15").Value = 16
$ws.Range("K15").Value = 'Mankaki'
$ws.Range("L15").Value = 'Especial'
$ws.Range("Q15").Value = '$/caja 16 kilos granel'
$ws.Range("R15").Value = 'Región Metropolitana'
$ws.Range("D16").Value = 45083
$ws.Range("M16").Value = 250
$ws.Range("N16").Value = 11200
$ws.Range("O16").Value = 11200
$ws.Range("P16").Value = 11200
$ws.Range("S16").Value = 700
$ws.Range("T16").Value = 16
$ws.Range("K16").Value = 'Mankaki'
$ws.Range("L16").Value = 'Primera'
$ws.Range("Q16").Value = '$/caja 16 kilos granel'
$ws.Range("R16").Value = 'Región Metropolitana'
$ws.Range("D17").Value = 45083
$ws.Range("M17").Value = 280
$ws.Range("N17").Value = 9600
$ws.Range("O17").Value = 9600
$ws.Range("P17").Value = 9600
$ws.Range("S17").Value = 600
$ws.Range("T17").Value = 16
$ws.Range("K17").Value = 'Mankaki'
$ws.Range("L17").Value = 'Segunda'
$ws.Range("Q17").Value = '$/caja 16 kilos granel'
$ws.Range("R17").Value = 'Región Metropolitana'
$ws.Range("D18").Value = 45049
$ws.Range("M18").Value = 500
$ws.Range("N18").Value = 13000
$ws.Range("O18").Value = 14000
$ws.Range("P18").Value = 13560
$ws.Range("S18").Value = 848
$ws.Range("T18").Value = 16
$ws.Range("K18").Value = 'Mankaki'
$ws.Range("L18").Value = 'Primera'
$ws.Range("Q18").Value = '$/caja 16 kilos granel'
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("D19").Value = 44329
$ws.Range("M19").Value = 3
$ws.Range("N19").Value = 250000
$ws.Range("O19").Value = 250000
$ws.Range("P19").Value = 250000
$ws.Range("S19").Value = 556
$ws.Range("T19").Value = 450
$ws.Range("K19").Value = 'Hachiya'
$ws.Range("L19").Value = 'Especial'
$ws.Range("Q19").Value = '$/bins (450 kilos)'
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("D20").Value = 44329
$ws.Range("M20").Value = 6
$ws.Range("N20").Value = 230000
$ws.Range("O20").Value = 230000
$ws.Range("P20").Value = 230000
$ws.Range("S20").Value = 511
$ws.Range("T20").Value = 450
$ws.Range("K20").Value = 'Hachiya'
$ws.Range("L20").Value = 'Primera'
$ws.Range("Q20").Value = '$/bins (450 kilos)'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("D21").Value = 44329
$ws.Range("M21").Value = 8
$ws.Range("N21").Value = 200000
$ws.Range("O21").Value = 200000
$ws.Range("P21").Value = 200000
$ws.Range("S21").Value = 444
$ws.Range("T21").Value = 450
$ws.Range("K21").Value = 'Hachiya'
$ws.Range("L21").Value = 'Segunda'
$ws.Range("Q21").Value = '$/bins (450 kilos)'
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("D22").Value = 44329
$ws.Range("M22").Value = 85
$ws.Range("N22").Value = 15000
$ws.Range("O22").Value = 15000
$ws.Range("P22").Value = 15000
$ws.Range("S22").Value = 1000
$ws.Range("T22").Value = 15
$ws.Range("K22").Value = 'Mankaki'
$ws.Range("L22").Value = 'Primera'
$ws.Range("Q22").Value = '$/caja 15 kilos granel'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("D23").Value = 44329
$ws.Range("M23").Value = 110
$ws.Range("N23").Value = 13000
$ws.Range("O23").Value = 13000
$ws.Range("P23").Value = 13000
$ws.Range("S23").Value = 867
$ws.Range("T23").Value = 15
$ws.Range("K23").Value = 'Mankaki'
$ws.Range("L23").Value = 'Segunda'
$ws.Range("Q23").Value = '$/caja 15 kilos granel'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("D24").Value = 44335
$ws.Range("M24").Value = 40
$ws.Range("N24").Value = 18000
$ws.Range("O24").Value = 18000
$ws.Range("P24").Value = 18000
$ws.Range("S24").Value = 1000
$ws.Range("T24").Value = 18
$ws.Range("K24").Value = 'Hachiya'
$ws.Range("L24").Value = 'Primera'
$ws.Range("Q24").Value = '$/caja 18 kilos granel'
$ws.Range("R24").Value = 'Región Metropolitana'
$ws.Range("D25").Value = 44335
$ws.Range("M25").Value = 55
$ws.Range("N25").Value = 14000
$ws.Range("O25").Value = 14000
$ws.Range("P25").Value = 14000
$ws.Range("S25").Value = 778
$ws.Range("T25").Value = 18
$ws.Range("K25").Value = 'Hachiya'
$ws.Range("L25").Value = 'Segunda'
$ws.Range("Q25").Value = '$/caja 18 kilos granel'
$ws.Range("R25").Value = 'Región Metropolitana'
$ws.Range("D26").Value = 44336
$ws.Range("M26").Value = 45
$ws.Range("N26").Value = 18000
$ws.Range("O26").Value = 18000
$ws.Range("P26").Value = 18000
$ws.Range("S26").Value = 1000
$ws.Range("T26").Value = 18
$ws.Range("K26").Value = 'Mankaki'
$ws.Range("L26").Value = 'Primera'
$ws.Range("Q26").Value = '$/caja 18 kilos granel'
$ws.Range("R26").Value = 'Región de O''Higgins'
$ws.Range("D27").Value = 44336
$ws.Range("M27").Value = 70
$ws.Range("N27").Value = 14000
$ws.Range("O27").Value = 14000
$ws.Range("P27").Value = 14000
$ws.Range("S27").Value = 778
$ws.Range("T27").Value = 18
$ws.Range("K27").Value = 'Mankaki'
$ws.Range("L27").Value = 'Segunda'
$ws.Range("Q27").Value = '$/caja 18 kilos granel'
$ws.Range("R27").Value = 'Región de O''Higgins'
$ws.Range("D28").Value = 45090
$ws.Range("M28").Value = 470
$ws.Range("N28").Value = 10500
$ws.Range("O28").Value = 11000
$ws.Range("P28").Value = 10734
$ws.Range("S28").Value = 716
$ws.Range("T28").Value = 15
$ws.Range("K28").Value = 'Fuyu'
$ws.Range("L28").Value = 'Primera'
$ws.Range("Q28").Value = '$/caja 15 kilos granel'
$ws.Range("R28").Value = 'Provincia de Curicó'
$ws.Range("D29").Value = 44328
$ws.Range("M29").Value = 55
$ws.Range("N29").Value = 15000
$ws.Range("O29").Value = 15000
$ws.Range("P29").Value = 15000
$ws.Range("S29").Value = 1000
$ws.Range("T29").Value = 15
$ws.Range("K29").Value = 'Mankaki'
$ws.Range("L29").Value = 'Primera'
$ws.Range("Q29").Value = '$/caja 15 kilos granel'
$ws.Range("R29").Value = 'Región de O''Higgins'
$ws.Range("D30").Value = 44328
$ws.Range("M30").Value = 80
$ws.Range("N30").Value = 13000
$ws.Range("O30").Value = 13000
$ws.Range("P30").Value = 13000
$ws.Range("S30").Value = 867
$ws.Range("T30").Value = 15
$ws.Range("K30").Value = 'Mankaki'
$ws.Range("L30").Value = 'Segunda'
$ws.Range("Q30").Value = '$/caja 15 kilos granel'
$ws.Range("R30").Value = 'Región de O''Higgins'
$ws.Range("D31").Value = 44330
$ws.Range("M31").Value = 55
$ws.Range("N31").Value = 14000
$ws.Range("O31").Value = 14000
$ws.Range("P31").Value = 14000
$ws.Range("S31").Value = 933
$ws.Range("T31").Value = 15
$ws.Range("K31").Value = 'Mankaki'
$ws.Range("L31").Value = 'Primera'
$ws.Range("Q31").Value = '$/caja 15 kilos granel'
$ws.Range("R31").Value = 'Región Metropolitana'
$ws.Range("D32").Value = 44330
$ws.Range("M32").Value = 70
$ws.Range("N32").Value = 12000
$ws.Range("O32").Value = 12000
$ws.Range("P32").Value = 12000
$ws.Range("S32").Value = 800
$ws.Range("T32").Value = 15
$ws.Range("K32").Value = 'Mankaki'
$ws.Range("L32").Value = 'Segunda'
$ws.Range("Q32").Value = '$/caja 15 kilos granel'
$ws.Range("R32").Value = 'Región Metropolitana'
$ws.Range("D33").Value = 44357
$ws.Range("M33").Value = 200
$ws.Range("N33").Value = 14000
$ws.Range("O33").Value = 14000
$ws.Range("P33").Value = 14000
$ws.Range("S33").Value = 1000
$ws.Range("T33").Value = 14
$ws.Range("K33").Value = 'Mankaki'
$ws.Range("L33").Value = 'Primera'
$ws.Range("Q33").Value = '$/caja 14 kilos granel'
$ws.Range("R33").Value = 'Región de O''Higgins'
$ws.Range("D34").Value = 45086
$ws.Range("M34").Value = 220
$ws.Range("N34").Value = 10500
$ws.Range("O34").Value = 10500
$ws.Range("P34").Value = 10500
$ws.Range("S34").Value = 700
$ws.Range("T34").Value = 15
$ws.Range("K34").Value = 'Mankaki'
$ws.Range("L34").Value = 'Primera'
$ws.Range("Q34").Value = '$/caja 15 kilos granel'
$ws.Range("R34").Value = 'Región de O''Higgins'
$ws.Range("D35").Value = 45086
$ws.Range("M35").Value = 280
$ws.Range("N35").Value = 9000
$ws.Range("O35").Value = 9000
$ws.Range("P35").Value = 9000
$ws.Range("S35").Value = 600
$ws.Range("T35").Value = 15
$ws.Range("K35").Value = 'Mankaki'
$ws.Range("L35").Value = 'Segunda'
$ws.Range("Q35").Value = '$/caja 15 kilos granel'
$ws.Range("R35").Value = 'Región de O''Higgins'
$ws.Range("D36").Value = 44351
$ws.Range("M36").Value = 4
$ws.Range("N36").Value = 260000
$ws.Range("O36").Value = 260000
$ws.Range("P36").Value = 260000
$ws.Range("S36").Value = 578
$ws.Range("T36").Value = 450
$ws.Range("K36").Value = 'Mankaki'
$ws.Range("L36").Value = 'Primera'
$ws.Range("Q36").Value = '$/bins (450 kilos)'
$ws.Range("R36").Value = 'Región de O''Higgins'
$ws.Range("D37").Value = 44344
$ws.Range("M37").Value = 6
$ws.Range("N37").Value = 290000
$ws.Range("O37").Value = 290000
$ws.Range("P37").Value = 290000
$ws.Range("S37").Value = 644
$ws.Range("T37").Value = 450
$ws.Range("K37").Value = 'Mankaki'
$ws.Range("L37").Value = 'Primera'
$ws.Range("Q37").Value = '$/bins (450 kilos)'
$ws.Range("R37").Value = 'Región Metropolitana'
$ws.Range("D38").Value = 45055
$ws.Range("M38").Value = 300
$ws.Range("N38").Value = 12500
$ws.Range("O38").Value = 12500
$ws.Range("P38").Value = 12500
$ws.Range("S38").Value = 781
$ws.Range("T38").Value = 16
$ws.Range("K38").Value = 'Mankaki'
$ws.Range("L38").Value = 'Primera'
$ws.Range("Q38").Value = '$/caja 16 kilos granel'
$ws.Range("R38").Value = 'Región de O''Higgins'
$ws.Range("D39").Value = 45055
$ws.Range("M39").Value = 280
$ws.Range("N39").Value = 9500
$ws.Range("O39").Value = 9500
$ws.Range("P39").Value = 9500
$ws.Range("S39").Value = 594
$ws.Range("T39").Value = 16
$ws.Range("K39").Value = 'Mankaki'
$ws.Range("L39").Value = 'Segunda'
$ws.Range("Q39").Value = '$/caja 16 kilos granel'
$ws.Range("R39").Value = 'Región de O''Higgins'
$ws.Range("D40").Value = 44301
$ws.Range("M40").Value = 120
$ws.Range("N40").Value = 18000
$ws.Range("O40").Value = 18000
$ws.Range("P40").Value = 18000
$ws.Range("S40").Value = 1000
$ws.Range("T40").Value = 18
$ws.Range("K40").Value = 'Fuyu'
$ws.Range("L40").Value = 'Primera'
$ws.Range("Q40").Value = '$/caja 18 kilos granel'
$ws.Range("R40").Value = 'Región de O''Higgins'
$ws.Range("D41").Value = 45071
$ws.Range("M41").Value = 330
$ws.Range("N41").Value = 15000
$ws.Range("O41").Value = 15000
$ws.Range("P41").Value = 15000
$ws.Range("S41").Value = 938
$ws.Range("T41").Value = 16
$ws.Range("K41").Value = 'Mankaki'
$ws.Range("L41").Value = 'Especial'
$ws.Range("Q41").Value = '$/caja 16 kilos granel'
$ws.Range("R41").Value = 'Región de O''Higgins'
$ws.Range("D42").Value = 45071
$ws.Range("M42").Value = 280
$ws.Range("N42").Value = 13000
$ws.Range("O42").Value = 13000
$ws.Range("P42").Value = 13000
$ws.Range("S42").Value = 812
$ws.Range("T42").Value = 16
$ws.Range("K42").Value = 'Mankaki'
$ws.Range("L42").Value = 'Primera'
$ws.Range("Q42").Value = '$/caja 16 kilos granel'
$ws.Range("R42").Value = 'Región de O''Higgins'
$ws.Range("D43").Value = 44316
$ws.Range("M43").Value = 4
$ws.Range("N43").Value = 300000
$ws.Range("O43").Value = 300000
$ws.Range("P43").Value = 300000
$ws.Range("S43").Value = 667
$ws.Range("T43").Value = 450
$ws.Range("K43").Value = 'Fuyu'
$ws.Range("L43").Value = 'Primera'
$ws.Range("Q43").Value = '$/bins (450 kilos)'
$ws.Range("R43").Value = 'Región Metropolitana'
$ws.Range("D44").Value = 45092
$ws.Range("M44").Value = 480
$ws.Range("N44").Value = 10000
$ws.Range("O44").Value = 10500
$ws.Range("P44").Value = 10208
$ws.Range("S44").Value = 681
$ws.Range("T44").Value = 15
$ws.Range("K44").Value = 'Mankaki'
$ws.Range("L44").Value = 'Primera'
$ws.Range("Q44").Value = '$/caja 15 kilos granel'
$ws.Range("R44").Value = 'Región de O''Higgins'
$ws.Range("D45").Value = 44322
$ws.Range("M45").Value = 70
$ws.Range("N45").Value = 15000
$ws.Range("O45").Value = 15000
$ws.Range("P45").Value = 15000
$ws.Range("S45").Value = 1250
$ws.Range("T45").Value = 12
$ws.Range("K45").Value = 'Fuyu'
$ws.Range("L45").Value = 'Especial'
$ws.Range("Q45").Value = '$/caja 12 kilos empedrada'
$ws.Range("R45").Value = 'Región de O''Higgins'
$ws.Range("D46").Value = 44322
$ws.Range("M46").Value = 90
$ws.Range("N46").Value = 13000
$ws.Range("O46").Value = 13000
$ws.Range("P46").Value = 13000
$ws.Range("S46").Value = 1083
$ws.Range("T46").Value = 12
$ws.Range("K46").Value = 'Fuyu'
$ws.Range("L46").Value = 'Primera'
$ws.Range("Q46").Value = '$/caja 12 kilos empedrada'
$ws.Range("R46").Value = 'Región de O''Higgins'
$ws.Range("D47").Value = 44322
$ws.Range("M47").Value = 6
$ws.Range("N47").Value = 270000
$ws.Range("O47").Value = 270000
$ws.Range("P47").Value = 270000
$ws.Range("S47").Value = 600
$ws.Range("T47").Value = 450
$ws.Range("K47").Value = 'Mankaki'
$ws.Range("L47").Value = 'Primera'
$ws.Range("Q47").Value = '$/bins (450 kilos)'
$ws.Range("R47").Value = 'Región de O''Higgins'
$ws.Range("D48").Value = 44334
$ws.Range("M48").Value = 50
$ws.Range("N48").Value = 19000
$ws.Range("O48").Value = 19000
$ws.Range("P48").Value = 19000
$ws.Range("S48").Value = 1056
$ws.Range("T48").Value = 18
$ws.Range("K48").Value = 'Mankaki'
$ws.Range("L48").Value = 'Especial'
$ws.Range("Q48").Value = '$/caja 18 kilos granel'
$ws.Range("R48").Value = 'Región de O''Higgins'
$ws.Range("D49").Value = 44334
$ws.Range("M49").Value = 65
$ws.Range("N49").Value = 18000
$ws.Range("O49").Value = 18000
$ws.Range("P49").Value = 18000
$ws.Range("S49").Value = 1000
$ws.Range("T49").Value = 18
$ws.Range("K49").Value = 'Mankaki'
$ws.Range("L49").Value = 'Primera'
$ws.Range("Q49").Value = '$/caja 18 kilos granel'
$ws.Range("R49").Value = 'Región de O''Higgins'
$ws.Range("D50").Value = 44334
$ws.Range("M50").Value = 80
$ws.Range("N50").Value = 14000
$ws.Range("O50").Value = 14000
$ws.Range("P50").Value = 14000
$ws.Range("S50").Value = 778
$ws.Range("T50").Value = 18
$ws.Range("K50").Value = 'Mankaki'
$ws.Range("L50").Value = 'Segunda'
$ws.Range("Q50").Value = '$/caja 18 kilos granel'
$ws.Range("R50").Value = 'Región de O''Higgins'
$ws.Range("D51").Value = 44323
$ws.Range("M51").Value = 70
$ws.Range("N51").Value = 15000
$ws.Range("O51").Value = 15000
$ws.Range("P51").Value = 15000
$ws.Range("S51").Value = 1000
$ws.Range("T51").Value = 15
$ws.Range("K51").Value = 'Mankaki'
$ws.Range("L51").Value = 'Primera'
$ws.Range("Q51").Value = '$/caja 15 kilos granel'
$ws.Range("R51").Value = 'Región de O''Higgins'
$ws.Range("D52").Value = 44326
$ws.Range("M52").Value = 50
$ws.Range("N52").Value = 16000
$ws.Range("O52").Value = 16000
$ws.Range("P52").Value = 16000
$ws.Range("S52").Value = 1067
$ws.Range("T52").Value = 15
$ws.Range("K52").Value = 'Mankaki'
$ws.Range("L52").Value = 'Primera'
$ws.Range("Q52").Value = '$/caja 15 kilos granel'
$ws.Range("R52").Value = 'Región de O''Higgins'
$ws.Range("D53").Value = 44314
$ws.Range("M53").Value = 4
$ws.Range("N53").Value = 310000
$ws.Range("O53").Value = 310000
$ws.Range("P53").Value = 310000
$ws.Range("S53").Value = 689
$ws.Range("T53").Value = 450
$ws.Range("K53").Value = 'Fuyu'
$ws.Range("L53").Value = 'Primera'
$ws.Range("Q53").Value = '$/bins (450 kilos)'
$ws.Range("R53").Value = 'Región de O''Higgins'
$ws.Range("D54").Value = 44314
$ws.Range("M54").Value = 5
$ws.Range("N54").Value = 320000
$ws.Range("O54").Value = 320000
$ws.Range("P54").Value = 320000
$ws.Range("S54").Value = 711
$ws.Range("T54").Value = 450
$ws.Range("K54").Value = 'Mankaki'
$ws.Range("L54").Value = 'Primera'
$ws.Range("Q54").Value = '$/bins (450 kilos)'
$ws.Range("R54").Value = 'Región de O''Higgins'
$ws.Range("D55").Value = 45089
$ws.Range("M55").Value = 150
$ws.Range("N55").Value = 13500
$ws.Range("O55").Value = 13500
$ws.Range("P55").Value = 13500
$ws.Range("S55").Value = 900
$ws.Range("T55").Value = 15
$ws.Range("K55").Value = 'Mankaki'
$ws.Range("L55").Value = 'Especial'
$ws.Range("Q55").Value = '$/caja 15 kilos granel'
$ws.Range("R55").Value = 'Región de O''Higgins'
$ws.Range("D56").Value = 45089
$ws.Range("M56").Value = 170
$ws.Range("N56").Value = 10500
$ws.Range("O56").Value = 10500
$ws.Range("P56").Value = 10500
$ws.Range("S56").Value = 700
$ws.Range("T56").Value = 15
$ws.Range("K56").Value = 'Mankaki'
$ws.Range("L56").Value = 'Primera'
$ws.Range("Q56").Value = '$/caja 15 kilos granel'
$ws.Range("R56").Value = 'Región de O''Higgins'
$ws.Range("D57").Value = 45089
$ws.Range("M57").Value = 180
$ws.Range("N57").Value = 9000
$ws.Range("O57").Value = 9000
$ws.Range("P57").Value = 9000
$ws.Range("S57").Value = 600
$ws.Range("T57").Value = 15
$ws.Range("K57").Value = 'Mankaki'
$ws.Range("L57").Value = 'Segunda'
$ws.Range("Q57").Value = '$/caja 15 kilos granel'
$ws.Range("R57").Value = 'Región de O''Higgins'
$ws.Range("D58").Value = 45050
$ws.Range("M58").Value = 400
$ws.Range("N58").Value = 9000
$ws.Range("O58").Value = 9500
$ws.Range("P58").Value = 9225
$ws.Range("S58").Value = 769
$ws.Range("T58").Value = 12
$ws.Range("K58").Value = 'Mankaki'
$ws.Range("L58").Value = 'Primera'
$ws.Range("Q58").Value = '$/caja 12 kilos empedrada'
$ws.Range("R58").Value = 'Región del Maule'
$ws.Range("D59").Value = 44309
$ws.Range("M59").Value = 60
$ws.Range("N59").Value = 14000
$ws.Range("O59").Value = 14000
$ws.Range("P59").Value = 14000
$ws.Range("S59").Value = 1167
$ws.Range("T59").Value = 12
$ws.Range("K59").Value = 'Mankaki'
$ws.Range("L59").Value = 'Especial'
$ws.Range("Q59").Value = '$/caja 12 kilos empedrada'
$ws.Range("R59").Value = 'Región de O''Higgins'
$ws.Range("D60").Value = 44309
$ws.Range("M60").Value = 80
$ws.Range("N60").Value = 12000
$ws.Range("O60").Value = 12000
$ws.Range("P60").Value = 12000
$ws.Range("S60").Value = 1000
$ws.Range("T60").Value = 12
$ws.Range("K60").Value = 'Mankaki'
$ws.Range("L60").Value = 'Primera'
$ws.Range("Q60").Value = '$/caja 12 kilos empedrada'
$ws.Range("R60").Value = 'Región de O''Higgins'
$ws.Range("D61").Value = 44698
$ws.Range("M61").Value = 180
$ws.Range("N61").Value = 22400
$ws.Range("O61").Value = 22400
$ws.Range("P61").Value = 22400
$ws.Range("S61").Value = 1400
$ws.Range("T61").Value = 16
$ws.Range("K61").Value = 'Mankaki'
$ws.Range("L61").Value = 'Especial'
$ws.Range("Q61").Value = '$/caja 16 kilos granel'
$ws.Range("R61").Value = 'Provincia de Curicó'
$ws.Range("D62").Value = 44698
$ws.Range("M62").Value = 150
$ws.Range("N62").Value = 19200
$ws.Range("O62").Value = 19200
$ws.Range("P62").Value = 19200
$ws.Range("S62").Value = 1200
$ws.Range("T62").Value = 16
$ws.Range("K62").Value = 'Mankaki'
$ws.Range("L62").Value = 'Primera'
$ws.Range("Q62").Value = '$/caja 16 kilos granel'
$ws.Range("R62").Value = 'Provincia de Curicó'
$ws.Range("D63").Value = 44698
$ws.Range("M63").Value = 200
$ws.Range("N63").Value = 16000
$ws.Range("O63").Value = 16000
$ws.Range("P63").Value = 16000
$ws.Range("S63").Value = 1000
$ws.Range("T63").Value = 16
$ws.Range("K63").Value = 'Mankaki'
$ws.Range("L63").Value = 'Segunda'
$ws.Range("Q63").Value = '$/caja 16 kilos granel'
$ws.Range("R63").Value = 'Provincia de Curicó'
$ws.Range("D64").Value = 45068
$ws.Range("M64").Value = 280
$ws.Range("N64").Value = 12000
$ws.Range("O64").Value = 12000
$ws.Range("P64").Value = 12000
$ws.Range("S64").Value = 800
$ws.Range("T64").Value = 15
$ws.Range("K64").Value = 'Fuyu'
$ws.Range("L64").Value = 'Especial'
$ws.Range("Q64").Value = '$/caja 15 kilos granel'
$ws.Range("R64").Value = 'Región de O''Higgins'
$ws.Range("D65").Value = 45068
$ws.Range("M65").Value = 350
$ws.Range("N65").Value = 9000
$ws.Range("O65").Value = 9000
$ws.Range("P65").Value = 9000
$ws.Range("S65").Value = 600
$ws.Range("T65").Value = 15
$ws.Range("K65").Value = 'Fuyu'
$ws.Range("L65").Value = 'Primera'
$ws.Range("Q65").Value = '$/caja 15 kilos granel'
$ws.Range("R65").Value = 'Región de O''Higgins'
$ws.Range("D66").Value = 44707
$ws.Range("M66").Value = 220
$ws.Range("N66").Value = 16000
$ws.Range("O66").Value = 16000
$ws.Range("P66").Value = 16000
$ws.Range("S66").Value = 1000
$ws.Range("T66").Value = 16
$ws.Range("K66").Value = 'Mankaki'
$ws.Range("L66").Value = 'Especial'
$ws.Range("Q66").Value = '$/caja 16 kilos granel'
$ws.Range("R66").Value = 'Provincia de Curicó'
$ws.Range("D67").Value = 44707
$ws.Range("M67").Value = 280
$ws.Range("N67").Value = 12800
$ws.Range("O67").Value = 12800
$ws.Range("P67").Value = 12800
$ws.Range("S67").Value = 800
$ws.Range("T67").Value = 16
$ws.Range("K67").Value = 'Mankaki'
$ws.Range("L67").Value = 'Primera'
$ws.Range("Q67").Value = '$/caja 16 kilos granel'
$ws.Range("R67").Value = 'Provincia de Curicó'
$ws.Range("D68").Value = 44707
$ws.Range("M68").Value = 250
$ws.Range("N68").Value = 11200
$ws.Range("O68").Value = 11200
$ws.Range("P68").Value = 11200
$ws.Range("S68").Value = 700
$ws.Range("T68").Value = 16
$ws.Range("K68").Value = 'Mankaki'
$ws.Range("L68").Value = 'Segunda'
$ws.Range("Q68").Value = '$/caja 16 kilos granel'
$ws.Range("R68").Value = 'Provincia de Curicó'
$ws.Range("D69").Value = 44333
$ws.Range("M69").Value = 50
$ws.Range("N69").Value = 19500
$ws.Range("O69").Value = 19500
$ws.Range("P69").Value = 19500
$ws.Range("S69").Value = 1083
$ws.Range("T69").Value = 18
$ws.Range("K69").Value = 'Mankaki'
$ws.Range("L69").Value = 'Especial'
$ws.Range("Q69").Value = '$/caja 18 kilos granel'
$ws.Range("R69").Value = 'Región de O''Higgins'
$ws.Range("D70").Value = 44333
$ws.Range("M70").Value = 85
$ws.Range("N70").Value = 18000
$ws.Range("O70").Value = 18000
$ws.Range("P70").Value = 18000
$ws.Range("S70").Value = 1000
$ws.Range("T70").Value = 18
$ws.Range("K70").Value = 'Mankaki'
$ws.Range("L70").Value = 'Primera'
$ws.Range("Q70").Value = '$/caja 18 kilos granel'
$ws.Range("R70").Value = 'Región de O''Higgins'
$ws.Range("D71").Value = 44333
$ws.Range("M71").Value = 100
$ws.Range("N71").Value = 14000
$ws.Range("O71").Value = 14000
$ws.Range("P71").Value = 14000
$ws.Range("S71").Value = 778
$ws.Range("T71").Value = 18
$ws.Range("K71").Value = 'Mankaki'
$ws.Range("L71").Value = 'Segunda'
$ws.Range("Q71").Value = '$/caja 18 kilos granel'
$ws.Range("R71").Value = 'Región de O''Higgins'
$ws.Range("D72").Value = 45063
$ws.Range("M72").Value = 420
$ws.Range("N72").Value = 13000
$ws.Range("O72").Value = 14000
$ws.Range("P72").Value = 13476
$ws.Range("S72").Value = 842
$ws.Range("T72").Value = 16
$ws.Range("K72").Value = 'Mankaki'
$ws.Range("L72").Value = 'Primera'
$ws.Range("Q72").Value = '$/caja 16 kilos granel'
$ws.Range("R72").Value = 'Región de O''Higgins'
$ws.Range("D73").Value = 44315
$ws.Range("M73").Value = 6
$ws.Range("N73").Value = 300000
$ws.Range("O73").Value = 300000
$ws.Range("P73").Value = 300000
$ws.Range("S73").Value = 667
$ws.Range("T73").Value = 450
$ws.Range("K73").Value = 'Fuyu'
$ws.Range("L73").Value = 'Primera'
$ws.Range("Q73").Value = '$/bins (450 kilos)'
$ws.Range("R73").Value = 'Región de O''Higgins'
$ws.Range("D74").Value = 44315
$ws.Range("M74").Value = 8
$ws.Range("N74").Value = 310000
$ws.Range("O74").Value = 310000
$ws.Range("P74").Value = 310000
$ws.Range("S74").Value = 689
$ws.Range("T74").Value = 450
$ws.Range("K74").Value = 'Mankaki'
$ws.Range("L74").Value = 'Primera'
$ws.Range("Q74").Value = '$/bins (450 kilos)'
$ws.Range("R74").Value = 'Región de O''Higgins'
$ws.Range("D75").Value = 45076
$ws.Range("M75").Value = 560
$ws.Range("N75").Value = 10500
$ws.Range("O75").Value = 11000
$ws.Range("P75").Value = 10750
$ws.Range("S75").Value = 717
$ws.Range("T75").Value = 15
$ws.Range("K75").Value = 'Fuyu'
$ws.Range("L75").Value = 'Primera'
$ws.Range("Q75").Value = '$/caja 15 kilos granel'
$ws.Range("R75").Value = 'Región Metropolitana'
$ws.Range("D76").Value = 45076
$ws.Range("M76").Value = 500
$ws.Range("N76").Value = 10500
$ws.Range("O76").Value = 11000
$ws.Range("P76").Value = 10780
$ws.Range("S76").Value = 719
$ws.Range("T76").Value = 15
$ws.Range("K76").Value = 'Hachiya'
$ws.Range("L76").Value = 'Primera'
$ws.Range("Q76").Value = '$/caja 15 kilos granel'
$ws.Range("R76").Value = 'Región de O''Higgins'
$ws.Range("D77").Value = 44327
$ws.Range("M77").Value = 8
$ws.Range("N77").Value = 260000
$ws.Range("O77").Value = 260000
$ws.Range("P77").Value = 260000
$ws.Range("S77").Value = 578
$ws.Range("T77").Value = 450
$ws.Range("K77").Value = 'Mankaki'
$ws.Range("L77").Value = 'Primera'
$ws.Range("Q77").Value = '$/bins (450 kilos)'
$ws.Range("R77").Value = 'Región de O''Higgins'
$ws.Range("D78").Value = 44327
$ws.Range("M78").Value = 80
$ws.Range("N78").Value = 15000
$ws.Range("O78").Value = 15000
$ws.Range("P78").Value = 15000
$ws.Range("S78").Value = 1000
$ws.Range("T78").Value = 15
$ws.Range("K78").Value = 'Mankaki'
$ws.Range("L78").Value = 'Primera'
$ws.Range("Q78").Value = '$/caja 15 kilos granel'
$ws.Range("R78").Value = 'Región de O''Higgins'
$ws.Range("D79").Value = 44327
$ws.Range("M79").Value = 7
$ws.Range("N79").Value = 220000
$ws.Range("O79").Value = 220000
$ws.Range("P79").Value = 220000
$ws.Range("S79").Value = 489
$ws.Range("T79").Value = 450
$ws.Range("K79").Value = 'Mankaki'
$ws.Range("L79").Value = 'Segunda'
$ws.Range("Q79").Value = '$/bins (450 kilos)'
$ws.Range("R79").Value = 'Región de O''Higgins'
$ws.Range("D80").Value = 44327
$ws.Range("M80").Value = 120
$ws.Range("N80").Value = 13000
$ws.Range("O80").Value = 13000
$ws.Range("P80").Value = 13000
$ws.Range("S80").Value = 867
$ws.Range("T80").Value = 15
$ws.Range("K80").Value = 'Mankaki'
$ws.Range("L80").Value = 'Segunda'
$ws.Range("Q80").Value = '$/caja 15 kilos granel'
$ws.Range("R80").Value = 'Región de O''Higgins'
$ws.Range("D81").Value = 45069
$ws.Range("M81").Value = 470
$ws.Range("N81").Value = 13000
$ws.Range("O81").Value = 13500
$ws.Range("P81").Value = 13234
$ws.Range("S81").Value = 827
$ws.Range("T81").Value = 16
$ws.Range("K81").Value = 'Fuyu'
$ws.Range("L81").Value = 'Primera'
$ws.Range("Q81").Value = '$/caja 16 kilos granel'
$ws.Range("R81").Value = 'Región de O''Higgins'
